$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header comment label to the longer description
$ws.Range("F1").Value = "Description/Comment"

# Add an example/instruction row beneath the headers
$ws.Range("A2").Value = "Sec"
$ws.Range("B2").Value = "Sub"
$ws.Range("C2").Value = "Dis"
$ws.Range("D2").Value = "Dev"
$ws.Range("E2").Value = "Idx"
$ws.Range("F2").Value = "replace this row with real data before uploading. All name elements, except for the instance index, needs to be approved in the Area and Device structrure of the Naming Service before import. "
